$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '67.134.15'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +6.38%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.559.42'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +10.28%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '188.82'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +9.77%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '552.58'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.57%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.551.69'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +10.12%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.609'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.59%  '
$ws.Range('E9').Value = '  -0.01%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.634'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +4.69%  '
$ws.Range('E11').Value = '  +14.05%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '54.69'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.82%  '
$ws.Range('E13').Value = '  +6.83%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.39'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.02%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.132.44'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +10.46%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.564.15'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +10.55%  '
$ws.Range('E17').Value = '  +4.95%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '67.140.85'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +6.54%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '18.23'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +6.15%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.00'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +8.58%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.995'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.86%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '434.30'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +18.76%  '
$ws.Range('B23').Value = 'PancakeSwap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.96'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +4.87%  '
$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '85.34'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +4.95%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.12'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.74%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.11'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.71%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.91'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +10.14%  '
$ws.Range('E28').Value = '  -0.54%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '12.15'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +7.42%  '
$ws.Range('E30').Value = '  +11.02%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '30.40'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +6.69%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '649.33'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.93%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.62'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.35%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '11.73'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +4.30%  '
$ws.Range('E35').Value = '  +5.52%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '59.64'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +4.91%  '
$ws.Range('B37').Value = 'PEPE'
$ws.Range('C37').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0₃0826'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +16.13%  '
$ws.Range('B38').Value = 'Kaspa'
$ws.Range('C38').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.150'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +23.03%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '38.70'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +5.23%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.999'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.11%  '
$ws.Range('E41').Value = '  +4.22%  '
$ws.Range('E42').Value = '  +14.64%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.00'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.05%  '
$ws.Range('E44').Value = '  +4.79%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.035.14'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +5.52%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.89'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +11.72%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.40'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +12.23%  '
$ws.Range('E48').Value = '  +6.77%  '
$ws.Range('E49').Value = '  +4.27%  '
$ws.Range('E50').Value = '  +5.17%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '141.93'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +5.29%  '
